$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the "Level 1: Location
#    Transparency" paragraph to the end of the "Abstract" paragraph (the
#    very first paragraph in the document).
# ---------------------------------------------------------------------------

# Remove the bookmark from its old location (end of the "Level 1: Location
# Transparency" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create it, collapsed, right after the text of the "Abstract" paragraph
# (i.e. after the run, before the paragraph mark) -- the same shape it used
# to have in the "Level 1..." paragraph. A collapsed range positioned
# exactly on the paragraph mark can't be targeted directly, so a
# one-character placeholder is appended to the run, a bookmark is wrapped
# around that placeholder (yielding a range collapsed to just after the
# real text once the placeholder is removed), and the placeholder text is
# then deleted again -- leaving the now-empty bookmark sitting in the right
# spot.
$abstractPara = $d.Paragraphs(1)
$endRange = $abstractPara.Range.Duplicate
$endRange.MoveEnd(1, -1)      # wdCharacter: exclude the paragraph mark
$endRange.Collapse(0)          # wdCollapseEnd: collapse to right after "Abstract"
$endRange.InsertAfter("~")
$endRange.Bookmarks.Add("_GoBack") | Out-Null
$endRange.Text = ""

# ---------------------------------------------------------------------------
# 2) Append the "Contribution" and "Conclusion" sections at the end of the
#    document (after the last, already-empty paragraph that currently
#    precedes </w:body>).
# ---------------------------------------------------------------------------

function Add-BodyParagraph($text) {
    $n = $d.Paragraphs.Count
    $d.Paragraphs($n).Range.InsertParagraphAfter()
    $n2 = $d.Paragraphs.Count
    $p = $d.Paragraphs($n2)
    if ($text -ne "") {
        $p.Range.InsertAfter($text)
    } else {
        # Leave a genuinely run-less empty paragraph, matching the style of
        # the document's other blank paragraphs: insert a placeholder
        # character and remove it again rather than leaving an untouched
        # (but COM-materialised) empty run behind.
        $p.Range.InsertAfter("~")
        $p2 = $d.Paragraphs($n2)
        $rr = $d.Range($p2.Range.Start, $p2.Range.Start + 1)
        $rr.Text = ""
    }
}

Add-BodyParagraph "Contribution"
Add-BodyParagraph "My part in the project was creating some demo Functions, Procedures and Triggers, though we did everything together. In my part I implemented above mentioned things and explained how these things work to my groupmates so that they understand what I did. And in the same way they have done their responsibilities and explained what they have tried to achieve and how those things work. At the same time we kept the report quiet similar and wrote our parts in the report individually."
Add-BodyParagraph ""
Add-BodyParagraph "Conclusion"
Add-BodyParagraph "The project we tried was based on a real life scenario of how largest online marketplace works. Though we couldn’t implement it as a whole but tried to give an idea how distributed database work in real life scenario. We look forward to implement the project in future on a larger scale."
